# Add two new columns (AV, AW) to Sheet1:
#   AV1 = "Diad1_PDF_Model", AW1 = "Diad2_PDF_Model"
#   AV2:AV17 = "PseudoVoigtModel", AW2:AW17 = "PseudoVoigtModel"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row (copy style from existing header cell AU1 so formatting matches)
$ws.Range("AV1").Value = "Diad1_PDF_Model"
$ws.Range("AW1").Value = "Diad2_PDF_Model"
$ws.Range("AU1").Copy()
$ws.Range("AV1:AW1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-17
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 48).Value = "PseudoVoigtModel"  # column AV = 48
    $ws.Cells.Item($r, 49).Value = "PseudoVoigtModel"  # column AW = 49
}
